$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STEPS")

# Insert a new column before column G (shifts TC_STEP_ACTION etc. one column right)
$ws.Columns("G").Insert()

# Set the header text for the newly inserted column
$ws.Range("G1").Value = "TC_STEP_CALL_DATASET"

# New column inherits the width of the column to its left (F)
$ws.Columns("G").ColumnWidth = $ws.Columns("F").ColumnWidth

# Select G2 to match resulting workbook selection state
$ws.Activate()
$ws.Range("G2").Select()
